$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Fix the "Números Reales" -> "Números reales" label across column A ---
# First normalize the one row whose font differs (row 67 used a slightly
# different font) by copying the format from a row that already uses the
# common style, so the whole column ends up visually consistent.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A67").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A2:A110").Value = "Números reales"
$ws.Range("A2:A110").HorizontalAlignment = -4131

# --- Highlight cell C5 in yellow (reviewer flagged this text as unnecessary) ---
$ws.Range("C5").Interior.Color = 65535

# --- Reviewer comments (Lzambrano) ---
$ws.Range("C5").AddComment("Lzambrano:`nEste texto sobra") | Out-Null
$ws.Range("C65").AddComment("Lzambrano:`nEste destacado y el recuerda no se pueden regresar a la sección 1. Son componentes de sección 2. Agregar texto antes del destacado para vincular estas secciones en el manscrito.") | Out-Null
$ws.Range("D86").AddComment("Lzambrano:`nIncluir línea de texto  ajustar guion") | Out-Null
$ws.Range("B87").AddComment("Lzambrano:`nColocar en el manuscrito una o varias líneas de texto antes del destacado.") | Out-Null
$ws.Range("E91").AddComment("Lzambrano:`nAcá que hay") | Out-Null
$ws.Range("F96").AddComment("Lzambrano:`nDebe quedar dentro de una sección 2. ¿Cuál es?") | Out-Null
$ws.Range("D106").AddComment("Lzambrano:`nIncluir línea de texto en manuscrito y registrarla en el esqueleto de guion") | Out-Null
$ws.Range("B108").AddComment("Lzambrano:`nCambiar nombre por Competencias") | Out-Null
$ws.Range("A109").AddComment("Lzambrano:`nEsto corresponde a fin de tema") | Out-Null

Write-Output "edit complete"
